# Update DateBase/orders/Fresh bloom Flowers_2024-11-16.xlsx
# - Append rows 32-38 of new flower order lines to the "Orders" sheet.
# - Update the aggregated order-code string in G2 of the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Orders": append new order rows (columns C and F) ----
$ws = $wb.Worksheets.Item("Orders")

$newRows = @(
    @(32, "457_茴香花_lace flower yellow_undefined_1bunch", "15"),
    @(33, "454_蓝星花_tweedia blue_undefined_1bunch", "35"),
    @(34, "389_金合欢_mimosa_undefined_1bunch", "10"),
    @(35, "1_白洋桔梗_White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g", "0"),
    @(36, "490_米花 粉_rice flower pink_undefined_1bunch", "15"),
    @(37, "465_羽衣甘蓝_Brassica_Brassica oleracea var. acephala DC._1bunch", "10"),
    @(38, "558_油画小菊_Helenium_undefined_1bunch", "15")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $name = $row[1]
    $qty = $row[2]

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $name

    $fCell = $ws.Cells.Item($r, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $qty
}

# ---- Sheet "Summary": extend the order-code string in G2 ----
$ws2 = $wb.Worksheets.Item("Summary")
$gCell = $ws2.Range("G2")
$gCell.NumberFormat = "@"
$gCell.Value = "02020351010101010101010222201051010635531030150101535100151015"
